# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 2;  3 = 0;  4 = 2;  5 = 1;  6 = 0;  7 = 1;  8 = 0;  9 = 3;  10 = 3;
    11 = 2; 12 = 1; 13 = 0; 14 = 2; 15 = 0; 16 = 1; 17 = 0; 18 = 0; 19 = 1;
    20 = 0; 21 = 0; 22 = 1; 23 = 0; 24 = 0; 25 = 0; 26 = 3; 27 = 1; 28 = 2;
    29 = 2; 30 = 2; 31 = 0; 32 = 1; 33 = 2; 34 = 1; 35 = 3; 36 = 0; 37 = 2;
    38 = 1; 39 = 0; 40 = 0; 41 = 2; 42 = 0; 43 = 1; 44 = 2; 45 = 2; 46 = 0;
    47 = 0; 48 = 1; 49 = 1; 50 = 1; 51 = 1; 52 = 2; 53 = 2; 54 = 0; 55 = 1;
    56 = 2; 57 = 2; 58 = 1; 59 = 1; 60 = 1; 61 = 1; 62 = 1; 63 = 1; 64 = 0;
    65 = 1; 66 = 2; 67 = 1; 68 = 1; 69 = 0; 70 = 0; 71 = 1; 72 = 1; 73 = 1;
    74 = 2; 76 = 3; 77 = 1; 78 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
